$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 6888.778
$ws.Range("J32").Value = 12249.75
$ws.Range("L32").Value = 12249.75
$ws.Range("N32").Value = -12901.75

$ws.Range("H82").Value = 2250
$ws.Range("I82").Value = 2250
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 6750
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -6344
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 2250
$ws.Range("I85").Value = 2250
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 6750
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -5346
$ws.Range("N85").ClearContents()

$ws.Range("H86").Value = 9750.5
$ws.Range("I86").Value = 3000.6667
$ws.Range("J86").Value = 30000
$ws.Range("K86").Value = 3000.6667
$ws.Range("L86").Value = 30000
$ws.Range("M86").Value = -1877.6667
$ws.Range("N86").Value = -32246

$ws.Range("H89").Value = 9750.5
$ws.Range("I89").Value = 3000.6667
$ws.Range("J89").Value = 30000
$ws.Range("K89").Value = 15003.3335
$ws.Range("L89").Value = 150000
$ws.Range("M89").Value = -9387.333500000001
$ws.Range("N89").Value = -161232

$ws.Range("H112").Value = 3180.8572
$ws.Range("J112").Value = 3264.111
$ws.Range("L112").Value = 9792.332999999999
$ws.Range("N112").Value = -12008.333

$ws.Range("H138").Value = 4740.304
$ws.Range("I138").Value = 4409.65
$ws.Range("J138").Value = 6944.6665
$ws.Range("K138").Value = 13228.95
$ws.Range("L138").Value = 20833.9995
$ws.Range("M138").Value = -8088.949999999999
$ws.Range("N138").Value = -31113.9995

$ws.Range("H141").Value = 5789.7334
$ws.Range("J141").Value = 12970.714
$ws.Range("L141").Value = 38912.142
$ws.Range("N141").Value = -49272.142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3078.5334
$ws.Range("I61").Value = 2874.3447
$ws.Range("K61").Value = 2874.3447
$ws.Range("M61").Value = -2662.3447

$ws.Range("H63").Value = 3795
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 3795
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H97").Value = 1558.8182
$ws.Range("I97").Value = 1194.1111
$ws.Range("K97").Value = 1194.1111
$ws.Range("M97").Value = -698.1111000000001

$ws.Range("H107").Value = 60000
$ws.Range("J107").Value = 60000
$ws.Range("L107").Value = 60000
$ws.Range("N107").Value = -67680

$ws.Range("H109").Value = 70000
$ws.Range("J109").Value = 70000
$ws.Range("L109").Value = 70000
$ws.Range("N109").Value = -72774

$ws.Range("H115").Value = 21561.334
$ws.Range("I115").Value = 5000
$ws.Range("J115").Value = 29842
$ws.Range("K115").Value = 5000
$ws.Range("L115").Value = 29842
$ws.Range("M115").Value = -3433
$ws.Range("N115").Value = -32976

$ws.Range("H136").Value = 3078.5334
$ws.Range("I136").Value = 2874.3447
$ws.Range("K136").Value = 8623.034100000001
$ws.Range("M136").Value = -6073.034100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2080.0417
$ws.Range("I105").Value = 1892.2106
$ws.Range("J105").Value = 2793.8
$ws.Range("K105").Value = 1892.2106
$ws.Range("L105").Value = 2793.8
$ws.Range("M105").Value = -145.2106000000001
$ws.Range("N105").Value = -6287.8

$ws.Range("H134").Value = 3319.5
$ws.Range("I134").Value = 3411.652
$ws.Range("K134").Value = 10234.956
$ws.Range("M134").Value = -7699.956

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H7").Value = 87.38461
$ws.Range("I7").Value = 40.285713
$ws.Range("J7").Value = 142.33333
$ws.Range("K7").Value = 40.285713
$ws.Range("L7").Value = 142.33333
$ws.Range("M7").Value = 72.714287
$ws.Range("N7").Value = -368.33333

$ws.Range("H25").Value = 10010
$ws.Range("I25").Value = 10010
$ws.Range("K25").Value = 10010
$ws.Range("M25").Value = -9836

$ws.Range("H31").Value = 2298.611
$ws.Range("I31").Value = 2605.75
$ws.Range("J31").Value = 1684.3334
$ws.Range("K31").Value = 2605.75
$ws.Range("L31").Value = 1684.3334
$ws.Range("M31").Value = -2310.75
$ws.Range("N31").Value = -2274.3334

$ws.Range("H34").Value = 2298.611
$ws.Range("I34").Value = 2605.75
$ws.Range("J34").Value = 1684.3334
$ws.Range("K34").Value = 2605.75
$ws.Range("L34").Value = 1684.3334
$ws.Range("M34").Value = -2403.75
$ws.Range("N34").Value = -2088.3334

$ws.Range("H58").Value = 5223.154
$ws.Range("I58").Value = 3484.889
$ws.Range("K58").Value = 3484.889
$ws.Range("M58").Value = -3281.889

$ws.Range("H96").Value = 11578.429
$ws.Range("J96").Value = 11578.429
$ws.Range("L96").Value = 11578.429
$ws.Range("N96").Value = -17070.429

$ws.Range("H132").Value = 4355.143
$ws.Range("J132").Value = 4257
$ws.Range("L132").Value = 12771
$ws.Range("N132").Value = -17831

$ws.Range("H134").Value = 3118.3635
$ws.Range("I134").Value = 3118.3635
$ws.Range("K134").Value = 9355.0905
$ws.Range("M134").Value = -6820.0905

$ws.Range("H136").Value = 5223.154
$ws.Range("I136").Value = 3484.889
$ws.Range("K136").Value = 10454.667
$ws.Range("M136").Value = -7904.667000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2714.5
$ws.Range("I107").Value = 93.59999999999999
$ws.Range("J107").Value = 3722.5386
$ws.Range("K107").Value = 280.8
$ws.Range("L107").Value = 11167.6158
$ws.Range("M107").Value = 1639.2
$ws.Range("N107").Value = -15007.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 4166721.8
$ws.Range("I2").Value = 5263195.5
$ws.Range("J2").Value = 122.2
$ws.Range("K2").Value = 5263195.5
$ws.Range("L2").Value = 122.2
$ws.Range("M2").Value = -5263082.5
$ws.Range("N2").Value = -348.2

$ws.Range("H5").Value = 15000
$ws.Range("I5").Value = 15000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -14888
$ws.Range("N5").ClearContents()

$ws.Range("H21").Value = 8331
$ws.Range("I21").Value = 4999
$ws.Range("K21").Value = 4999
$ws.Range("M21").Value = -4826

$ws.Range("H30").Value = 8331
$ws.Range("I30").Value = 4999
$ws.Range("K30").Value = 4999
$ws.Range("M30").Value = -4894

$ws.Range("H57").Value = 20113.428
$ws.Range("J57").Value = 20113.428
$ws.Range("L57").Value = 20113.428
$ws.Range("N57").Value = -21753.428

$ws.Range("H113").Value = 1611.4546
$ws.Range("I113").Value = 1580.6666
$ws.Range("K113").Value = 1580.6666
$ws.Range("M113").Value = 589.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 4006.5
$ws.Range("I32").Value = 4006.5
$ws.Range("K32").Value = 4006.5
$ws.Range("M32").Value = -3689.5

$ws.Range("H61").Value = 10252.143
$ws.Range("I61").Value = 9570.223
$ws.Range("K61").Value = 9570.223
$ws.Range("M61").Value = -9368.223

$ws.Range("H68").Value = 6923.952
$ws.Range("J68").Value = 7321
$ws.Range("L68").Value = 7321
$ws.Range("N68").Value = -8819

$ws.Range("H71").Value = 6923.952
$ws.Range("J71").Value = 7321
$ws.Range("L71").Value = 36605
$ws.Range("N71").Value = -44093

$ws.Range("H113").Value = 10252.143
$ws.Range("I113").Value = 9570.223
$ws.Range("K113").Value = 9570.223
$ws.Range("M113").Value = -7400.223

$ws.Range("H125").Value = 76837.5
$ws.Range("J125").Value = 76837.5
$ws.Range("L125").Value = 76837.5
$ws.Range("N125").Value = -86677.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 29930
$ws.Range("I2").Value = 29930
$ws.Range("K2").Value = 29930
$ws.Range("M2").Value = -29818

$ws.Range("H126").Value = 2094
$ws.Range("I126").Value = 1969.5333
$ws.Range("J126").Value = 2716.3333
$ws.Range("K126").Value = 5908.5999
$ws.Range("L126").Value = 8148.999899999999
$ws.Range("M126").Value = -3438.5999
$ws.Range("N126").Value = -13088.9999

$ws.Range("H132").Value = 1200.7715
$ws.Range("I132").Value = 936.73334
$ws.Range("K132").Value = 2810.20002
$ws.Range("M132").Value = -280.2000200000002

Write-Output "Applied all cell updates."
